$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates (Week 2/3/4: fill in PPG + FG% for the two Football
#     athletes, which were previously blank, and refresh the recalculated
#     Sprint Speed / FG% figures for the other rows) ---

# Athlete B (row 3) - add PPG and FG%
$ws.Range("C3").Value = 22.4
$ws.Range("D3").Value = 0.345
$ws.Range("E3").Value = 20.4

# Athlete A (row 2) - refreshed FG%
$ws.Range("D2").Value = 0.475

# Athlete C (row 4) - refreshed Sprint Speed
$ws.Range("E4").Value = 19.1

# Athlete D (row 5) - add PPG and FG%
$ws.Range("C5").Value = 27.4
$ws.Range("D5").Value = 0.535

# --- UI state: leave the selection on H6, matching the author's last
#     position when they saved the workbook ---
$ws.Range("H6").Select()
